$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of BOM data appended below the existing table (rows 52-57)
# Written in an order that matches the shared-string insertion order of the target file.
$ws.Range("C52").Value = "IC TOUCH SENSOR 1KEY SOT23-6"
$ws.Range("D52").Value = "AT42QT1011-TSHR"
$ws.Range("E52").Value = "AT42QT1011-TSHRCT-ND"
$ws.Range("F52").Value = "Atmel"
$ws.Range("B52").Value = "2"
$ws.Range("C52").Style = "Normal"

$ws.Range("C53").Value = "2.2nF Cap"
$ws.Range("B53").Value = "2"

$ws.Range("C54").Value = "22k RES"
$ws.Range("B54").Value = "2"

$ws.Range("C55").Value = "100nF Cap"
$ws.Range("B55").Value = "2"

$ws.Range("B56").Value = "1"
$ws.Range("C56").Value = "CONN HEADER PH SIDE 2POS 2MM"
$ws.Range("E56").Value = "455-1719-ND"
$ws.Range("C56").Style = "Normal"
$ws.Range("E56").Style = "Normal"

$ws.Range("B57").Value = "1"
$ws.Range("C57").Value = "IC MCU ARM 2MB FLASH 100LQFP"
$ws.Range("D57").Value = "STM32F429VIT6"
$ws.Range("E57").Value = "497-14052-ND"

# Update the view: scroll to show new rows, select C58
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("C58").Select()
